# Weekly update: insert a new daily price record as row 18 for
# "Vega Monumental Concepción" - Pepino ensalada, shifting the existing
# rows 18-95 down to 19-96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 18 through the end of the data down by one row.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new record's data.
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = "Vega Monumental Concepción"
$ws.Range("C18").Value = "Bíobío"
$ws.Range("D18").Value = 44558
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = 100112043
$ws.Range("G18").Value = "Pepino ensalada"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 250
$ws.Range("K18").Value = 8000
$ws.Range("L18").Value = 8500
$ws.Range("M18").Value = 8300
$ws.Range("N18").Value = "$/caja 60 unidades"
$ws.Range("O18").Value = "Provincia de Limarí"
$ws.Range("P18").Value = 138
$ws.Range("Q18").Value = 60
$ws.Range("R18").Value = "Hortaliza"
